$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.703179
$ws.Range("H2").Value = 2.109537
$ws.Range("I2").Value = 0.09061565978662672
$ws.Range("J2").Value = 0.09061565978662675
$ws.Range("M2").Value = 0.668273
$ws.Range("N2").Value = 2.004819
$ws.Range("O2").Value = 0.01328414746766746
$ws.Range("P2").Value = 0.01328414746766746
$ws.Range("Q2").Value = 0.469915539867
$ws.Range("R2").Value = 4.229239858803
$ws.Range("S2").Value = 0.001203751787485533
$ws.Range("T2").Value = 0.001203751787485533

# Row 3
$ws.Range("G3").Value = 0.703179
$ws.Range("H3").Value = 2.109537
$ws.Range("I3").Value = 0.09061565978662672
$ws.Range("J3").Value = 0.09061565978662675
$ws.Range("O3").Value = 0.3831531055114357
$ws.Range("P3").Value = 0.3831531055114357
$ws.Range("Q3").Value = 13.553718736287
$ws.Range("R3").Value = 121.983468626583
$ws.Range("S3").Value = 0.03471967145521375
$ws.Range("T3").Value = 0.03471967145521376

# Row 4
$ws.Range("G4").Value = 0.703179
$ws.Range("H4").Value = 2.109537
$ws.Range("I4").Value = 0.09061565978662672
$ws.Range("J4").Value = 0.09061565978662675
$ws.Range("M4").Value = 30.36285833333334
$ws.Range("N4").Value = 91.08857500000001
$ws.Range("O4").Value = 0.6035627470208969
$ws.Range("P4").Value = 0.6035627470208967
$ws.Range("Q4").Value = 21.350524359975
$ws.Range("R4").Value = 192.154719239775
$ws.Range("S4").Value = 0.05469223654392744
$ws.Range("T4").Value = 0.05469223654392745

# Row 5
$ws.Range("I5").Value = 0.6019120921953386
$ws.Range("J5").Value = 0.6019120921953387
$ws.Range("M5").Value = 0.668273
$ws.Range("N5").Value = 2.004819
$ws.Range("O5").Value = 0.01328414746766746
$ws.Range("P5").Value = 0.01328414746766746
$ws.Range("Q5").Value = 3.121401382746333
$ws.Range("R5").Value = 28.092612444717
$ws.Range("S5").Value = 0.007995888995295131
$ws.Range("T5").Value = 0.007995888995295129

# Row 6
$ws.Range("I6").Value = 0.6019120921953386
$ws.Range("J6").Value = 0.6019120921953387
$ws.Range("O6").Value = 0.3831531055114357
$ws.Range("P6").Value = 0.3831531055114357
$ws.Range("S6").Value = 0.2306244873695296
$ws.Range("T6").Value = 0.2306244873695296

# Row 7
$ws.Range("I7").Value = 0.6019120921953386
$ws.Range("J7").Value = 0.6019120921953387
$ws.Range("M7").Value = 30.36285833333334
$ws.Range("N7").Value = 91.08857500000001
$ws.Range("O7").Value = 0.6035627470208969
$ws.Range("P7").Value = 0.6035627470208967
$ws.Range("S7").Value = 0.3632917158305139
$ws.Range("T7").Value = 0.3632917158305139

# Row 8
$ws.Range("I8").Value = 0.3074722480180346
$ws.Range("J8").Value = 0.3074722480180347
$ws.Range("M8").Value = 0.668273
$ws.Range("N8").Value = 2.004819
$ws.Range("O8").Value = 0.01328414746766746
$ws.Range("P8").Value = 0.01328414746766746
$ws.Range("Q8").Value = 1.594492472512333
$ws.Range("R8").Value = 14.350432252611
$ws.Range("S8").Value = 0.004084506684886796
$ws.Range("T8").Value = 0.004084506684886796

# Row 9
$ws.Range("I9").Value = 0.3074722480180346
$ws.Range("J9").Value = 0.3074722480180347
$ws.Range("O9").Value = 0.3831531055114357
$ws.Range("P9").Value = 0.3831531055114357
$ws.Range("S9").Value = 0.1178089466866924
$ws.Range("T9").Value = 0.1178089466866924

# Row 10
$ws.Range("I10").Value = 0.3074722480180346
$ws.Range("J10").Value = 0.3074722480180347
$ws.Range("M10").Value = 30.36285833333334
$ws.Range("N10").Value = 91.08857500000001
$ws.Range("O10").Value = 0.6035627470208969
$ws.Range("P10").Value = 0.6035627470208967
$ws.Range("Q10").Value = 72.44546623379723
$ws.Range("R10").Value = 652.009196104175
$ws.Range("S10").Value = 0.1855787946464555
$ws.Range("T10").Value = 0.1855787946464555
